$d = $word.ActiveDocument

# --- 1. Text correction in the Type III compensator design paragraph ---
$old = "sawtooth wave has to be selected beforehand and which is selected as 1.8 Volts which is a standard selection. "
$new = "sawtooth wave and the reference voltage have to be selected beforehand and which are selected as 1.8 Volts and 1.2 Volts respectively which are standard selections. "

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) {
    throw "Could not find the target sentence to replace."
}

# --- 2. Resize Figure 3 (Type III compensator circuit schematic) ---
$shp = $d.InlineShapes.Item(3)
$shp.Width = 269.21732283464564
$shp.Height = 177.40330708661418
